$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-20 Friday" "2025-06-21 Saturday"

Replace-Text "122×4=" "323×8="
Replace-Text "385×5=" "372×4="
Replace-Text "766×6=" "730×7="
Replace-Text "140×7=" "977×9="
Replace-Text "810×6=" "910×4="

Replace-Text "757×8=" "760×4="
Replace-Text "608×8=" "680×9="
Replace-Text "308×2=" "960×6="
Replace-Text "991×3=" "579×2="
Replace-Text "883×9=" "256×8="

Replace-Text "951×2=" "630×9="
Replace-Text "244×2=" "966×6="
Replace-Text "501×7=" "299×6="
Replace-Text "828×9=" "828×5="
Replace-Text "969×2=" "170×5="

Replace-Text "329×4=" "735×6="
Replace-Text "141×3=" "579×6="
Replace-Text "751×8=" "259×6="
Replace-Text "771×2=" "727×8="
Replace-Text "808×5=" "337×7="

Replace-Text "754×3=" "526×4="
Replace-Text "179×6=" "901×7="
Replace-Text "443×8=" "966×6="
Replace-Text "916×7=" "671×3="
Replace-Text "657×8=" "780×6="

Write-Host "Done"
